$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.643.86"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "1.798.72"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5374"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3775"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.166"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.411"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "1.794.12"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06439"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.928"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "28.647.55"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.105"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.381"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "2.003.94"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.106"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1030"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.656"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.693"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("B35").Value = "Algorand"
$ws.Range("C35").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2258"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.63%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06491"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.878"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02315"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.034"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.210"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6246"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5869"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.957"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.158"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06893"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
